$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# constants_evaluated : Constant / St.Deviation / Validity  ->
#                        Component / Constant / St.Deviation / Validity
# collapse the 3 data rows into a single data row, add a "Component"
# name column, and store the constant/deviation as text (matches the
# regenerated regression output).
# ------------------------------------------------------------------
$wsConst = $wb.Worksheets.Item("constants_evaluated")
$wsConst.Range("A1:D4").ClearContents()
$wsConst.Range("A1").Value = "Component"
$wsConst.Range("B1").Value = "Constant"
$wsConst.Range("C1").Value = "St.Deviation"
$wsConst.Range("D1").Value = "Validity"

$wsConst.Range("A2").Value = "Comp"
$wsConst.Range("B2").NumberFormat = "@"
$wsConst.Range("B2").Value = "5.12441253662109"
$wsConst.Range("B2").Style = "Normal"
$wsConst.Range("C2").NumberFormat = "@"
$wsConst.Range("C2").Value = "0.0266629716825117"
$wsConst.Range("C2").Style = "Normal"
$wsConst.Range("D2").Value = "OK"

# ------------------------------------------------------------------
# enthalpies_calculated : reaction / value / dev
# collapse the 3 data rows into a single data row (PLP & T3H rows
# removed), and flip the sign of the Comp enthalpy value.
# ------------------------------------------------------------------
$wsEnth = $wb.Worksheets.Item("enthalpies_calculated")
$wsEnth.Range("A1:C4").ClearContents()
$wsEnth.Range("A1").Value = "reaction"
$wsEnth.Range("B1").Value = "value"
$wsEnth.Range("C1").Value = "dev"

$wsEnth.Range("A2").Value = "Comp"
$wsEnth.Range("B2").Value = 47.1839051361544
$wsEnth.Range("C2").Value = 0.587867672633848

# ------------------------------------------------------------------
# input_stoich_coefficients : PLP / T3H / name
# collapse the 3 data rows into a single data row; the remaining row
# stores "1"/"1" as text values rather than numeric 1/0.
# ------------------------------------------------------------------
$wsStoich = $wb.Worksheets.Item("input_stoich_coefficients")
$wsStoich.Range("A1:C4").ClearContents()
$wsStoich.Range("A1").Value = "PLP"
$wsStoich.Range("B1").Value = "T3H"
$wsStoich.Range("C1").Value = "name"

$wsStoich.Range("A2").NumberFormat = "@"
$wsStoich.Range("A2").Value = "1"
$wsStoich.Range("A2").Style = "Normal"
$wsStoich.Range("B2").NumberFormat = "@"
$wsStoich.Range("B2").Value = "1"
$wsStoich.Range("B2").Style = "Normal"
$wsStoich.Range("C2").Value = "Comp"

# ------------------------------------------------------------------
# input_enthalpies : reaction / value
# the PLP data row is cleared out (left blank).
# ------------------------------------------------------------------
$wsInEnth = $wb.Worksheets.Item("input_enthalpies")
$wsInEnth.Range("A2").ClearContents()
$wsInEnth.Range("B2").ClearContents()
